$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (country index 14)
$ws.Range("B2").Value2 = 0.232097580732653
$ws.Range("L2").Value2 = 0.386272295961198

# Row 3 (country index 15)
$ws.Range("B3").Value2 = 0.221547475344095
$ws.Range("L3").Value2 = 0.112246507177215

# Row 4 (country index 16)
$ws.Range("B4").Value2 = 0.318030317882592
$ws.Range("L4").Value2 = 0.22699232463484

# Row 5 (country index 17)
$ws.Range("B5").Value2 = 0.172370897143246
$ws.Range("E5").Value2 = 0.166158467826809
$ws.Range("L5").Value2 = 0.195595011071219
